# Docx writer: Use different style for block quotes in notes.
#
# Adds a new "Footnote Block Text" paragraph style (styleId
# "FootnoteBlockText"), based on "Footnote Text", with the same
# spacing/indent overrides already used by "Block Text", so footnote
# block quotes can be restyled independently of body block quotes.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$style = $d.Styles.Add("FootnoteBlockText", 1)

$style.NameLocal          = "Footnote Block Text"
$style.BaseStyle          = "Footnote Text"
$style.NextParagraphStyle = "Footnote Text"
$style.Priority           = 9
$style.UnhideWhenUsed     = $true
$style.QuickStyle         = $true

# Paragraph formatting is stored in twentieths of a point (twips) in the
# XML, but ParagraphFormat exposes it in points, so divide by 20.
$pf = $style.ParagraphFormat
$pf.SpaceBefore      = 100 / 20
$pf.SpaceAfter       = 100 / 20
$pf.FirstLineIndent  = 0 / 20
$pf.LeftIndent       = 480 / 20
$pf.RightIndent      = 480 / 20
